# Inserts one new weekly record as the most-recent row (row 371) in the
# "Hortaliza, Femacal de La Calera - Ajo" sheet, pushing the existing
# rows 371:459 down to 372:460.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 371:459 down to 372:460, leaving a blank row 371 behind.
$ws.Rows.Item(371).Insert()

# Populate the new row 371 with this week's record.
$ws.Cells.Item(371, 1).Value  = 3
$ws.Cells.Item(371, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(371, 3).Value  = "Coquimbo"
$ws.Cells.Item(371, 4).Value  = 44722
$ws.Cells.Item(371, 5).Value  = 5
$ws.Cells.Item(371, 6).Value  = 100112003
$ws.Cells.Item(371, 7).Value  = "Ajo"
$ws.Cells.Item(371, 8).Value  = "Chino"
$ws.Cells.Item(371, 9).Value  = "Primera"
$ws.Cells.Item(371, 10).Value = 82
$ws.Cells.Item(371, 11).Value = 17500
$ws.Cells.Item(371, 12).Value = 18000
$ws.Cells.Item(371, 13).Value = 17744
$ws.Cells.Item(371, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(371, 15).Value = "China"
$ws.Cells.Item(371, 16).Value = 1774
$ws.Cells.Item(371, 17).Value = 10
$ws.Cells.Item(371, 18).Value = "Hortaliza"
